$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '< -0.0898'
$ws.Range("E2").Value = '< -0.0893'
$ws.Range("F2").Value = '< -0.0898'
$ws.Range("G2").Value = '< -0.0797'
$ws.Range("H2").Value = '< -0.0802'
$ws.Range("I2").Value = '< -0.0797'
$ws.Range("J2").Value = '< -0.0888'
$ws.Range("M2").Value = '< -0.0888'
$ws.Range("O2").Value = '< -0.0888'
$ws.Range("G3").Value = '< -0.0264'
$ws.Range("I3").Value = '< -0.0264'
$ws.Range("J3").Value = '0.0879 -0.0355'
$ws.Range("M3").Value = '0.0879 -0.0355'
$ws.Range("O3").Value = '0.0879 -0.0355'
$ws.Range("E4").Value = '0.3612 5e-04'
$ws.Range("G4").Value = '0.5444 0.0101'
$ws.Range("H4").Value = '0.6255 0.0096'
$ws.Range("I4").Value = '0.5444 0.0101'
$ws.Range("J4").Value = '0.2554 0.001'
$ws.Range("K4").Value = '0.4527 0.0014'
$ws.Range("L4").Value = '0.4527 0.0014'
$ws.Range("M4").Value = '0.2554 0.001'
$ws.Range("N4").Value = '0.4527 0.0014'
$ws.Range("O4").Value = '0.2554 0.001'
$ws.Range("F5").Value = '0.6188 -5e-04'
$ws.Range("H5").Value = '0.7877 0.0091'
$ws.Range("J5").Value = '0.4909 5e-04'
$ws.Range("K5").Value = '0.2083 9e-04'
$ws.Range("L5").Value = '0.2083 9e-04'
$ws.Range("M5").Value = '0.4909 5e-04'
$ws.Range("N5").Value = '0.2083 9e-04'
$ws.Range("O5").Value = '0.4909 5e-04'
$ws.Range("G6").Value = '0.5444 0.0101'
$ws.Range("H6").Value = '0.6255 0.0096'
$ws.Range("I6").Value = '0.5444 0.0101'
$ws.Range("J6").Value = '0.2554 0.001'
$ws.Range("K6").Value = '0.4527 0.0014'
$ws.Range("L6").Value = '0.4527 0.0014'
$ws.Range("M6").Value = '0.2554 0.001'
$ws.Range("N6").Value = '0.4527 0.0014'
$ws.Range("O6").Value = '0.2554 0.001'
$ws.Range("H7").Value = '0.6139 -5e-04'
$ws.Range("J7").Value = '0.3545 -0.0091'
$ws.Range("K7").Value = '< -0.0086'
$ws.Range("L7").Value = '< -0.0086'
$ws.Range("M7").Value = '0.3545 -0.0091'
$ws.Range("N7").Value = '< -0.0086'
$ws.Range("O7").Value = '0.3545 -0.0091'
$ws.Range("I8").Value = '0.3659 5e-04'
$ws.Range("J8").Value = '0.2584 -0.0086'
$ws.Range("K8").Value = '0.1984 -0.0082'
$ws.Range("L8").Value = '0.1984 -0.0082'
$ws.Range("M8").Value = '0.2584 -0.0086'
$ws.Range("N8").Value = '0.1984 -0.0082'
$ws.Range("O8").Value = '0.2584 -0.0086'
$ws.Range("J9").Value = '0.3545 -0.0091'
$ws.Range("K9").Value = '< -0.0086'
$ws.Range("L9").Value = '< -0.0086'
$ws.Range("M9").Value = '0.3545 -0.0091'
$ws.Range("N9").Value = '< -0.0086'
$ws.Range("O9").Value = '0.3545 -0.0091'
$ws.Range("K10").Value = '0.5285 4e-04'
$ws.Range("L10").Value = '0.5285 4e-04'
$ws.Range("N10").Value = '0.5285 4e-04'
$ws.Range("M11").Value = '0.4515 -4e-04'
$ws.Range("O11").Value = '0.4515 -4e-04'
$ws.Range("M12").Value = '0.4515 -4e-04'
$ws.Range("O12").Value = '0.4515 -4e-04'
$ws.Range("N13").Value = '0.5285 4e-04'
$ws.Range("O14").Value = '0.4515 -4e-04'
